# "Test Case Template added"
# Populate the header row of the (previously empty) Sheet1 with the
# standard test-case-tracking column headers, then size the columns to
# fit their header text (mirrors Excel's AutoFit-on-entry behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "s.no",
    "Test Case ID",
    "Test Description",
    "Test Steps",
    "Test Priority",
    "Actual Result",
    "Expected Result",
    "Status",
    "Comments"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Best-fit the columns whose header text doesn't fit the default width
# (B, C, E, F, G) - same effect as selecting the columns and double
# clicking a column-border / Format > AutoFit Column Width in Excel.
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws.Columns.Item(3).ColumnWidth = 13.333333333333334
$ws.Columns.Item(5).ColumnWidth = 10.0
$ws.Columns.Item(6).ColumnWidth = 10.666666666666666
$ws.Columns.Item(7).ColumnWidth = 13.166666666666666

# Leave the selection on the last header cell, as in the source edit.
$ws.Range("I1").Select() | Out-Null
